$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$addresses = @(
    "1123 Hayden Meadows Drive,  OR 97217",
    "2201 Grand Blvd,  WA 98661",
    "17275 Nw Cornell Rd,  OR 97006",
    "9000 Ne Highway 99,  WA 98665",
    "7809 Ne Vancouver Plaza Dr,  WA 98662",
    "7650 Ne Shaleen Street,  OR 97006",
    "221e Ne 104th Ave,  WA 98664",
    "9055 Sw Murray Blvd,  OR 97008",
    "7600 Sw Dartmouth St.,  OR 97223",
    "4200 Se 82nd Ave,  OR 97266",
    "14505 Ne Fourth Plain Blvd,  WA 98682",
    "430 Se 192nd Ave,  WA 98683",
    "10000 Se 82nd Ave,  OR 97086",
    "15600 Se Mcloughlin Blvd,  OR 97267",
    "1201 Sw 13th Avenue,  WA 98604",
    "19133 Willamette Dr,  OR 97068",
    "220 N Adair St,  OR 97113",
    "3900 W Powell Blvd,  OR 97030",
    "21320 Sw Langer Farms Pkwy,  OR 97140",
    "23500 Ne Sandy Blvd,  OR 97060",
    "2295 Gable Rd,  OR 97051",
    "2444 E Powell Blvd,  OR 97080",
    "1486 Dike Access Rd,  WA 98674",
    "3002 Stacey Allison Way,  OR 97071",
    "2375 Ne Highway 99w,  OR 97128",
    "540 7th Ave,  WA 98632",
    "12620 Se 41st Pl,  WA 98006",
    "15063 Main St,  WA 98007",
    "743 Rainier Avenue South,  WA 98057",
    "6797 State Highway 303 Ne,  WA 98311",
    "3497 Bethel Rd Se,  WA 98366",
    "17222 Highway 99,  WA 98037",
    "1900 S 314th St,  WA 98003",
    "1400 164th St Sw,  WA 98087",
    "34520 16th Ave S,  WA 98003",
    "762 Outlet Collection Way,  WA 98001",
    "21200 Olhava Way Nw,  WA 98370",
    "11400 Highway 99,  WA 98204",
    "1605 Se Everett Mall Way,  WA 98208",
    "19191 N Kelsey Street,  WA 98272",
    "1965 S. Union Ave,  WA 98405",
    "7001 Bridgeport Way W,  WA 98499",
    "310 31st Ave Se,  WA 98374",
    "19205 State Route 410 E,  WA 98391",
    "16502 Meridian E,  WA 98375",
    "8713 64th St Ne,  WA 98270",
    "8924 Quilceda Blvd,  WA 98271",
    "4010 172nd St Ne,  WA 98223",
    "1401 Galaxy Dr Ne,  WA 98516",
    "5110 Yelm Highway,  WA 98503",
    "100 E Wallace Kneeland Blvd,  WA 98584",
)

$zips = @(
    "`"97217`"",
    "`"98661`"",
    "`"97006`"",
    "`"98665`"",
    "`"98662`"",
    "`"97006`"",
    "`"98664`"",
    "`"97008`"",
    "`"97223`"",
    "`"97266`"",
    "`"98682`"",
    "`"98683`"",
    "`"97086`"",
    "`"97267`"",
    "`"98604`"",
    "`"97068`"",
    "`"97113`"",
    "`"97030`"",
    "`"97140`"",
    "`"97060`"",
    "`"97051`"",
    "`"97080`"",
    "`"98674`"",
    "`"97071`"",
    "`"97128`"",
    "`"98632`"",
    "`"98006`"",
    "`"98007`"",
    "`"98057`"",
    "`"98311`"",
    "`"98366`"",
    "`"98037`"",
    "`"98003`"",
    "`"98087`"",
    "`"98003`"",
    "`"98001`"",
    "`"98370`"",
    "`"98204`"",
    "`"98208`"",
    "`"98272`"",
    "`"98405`"",
    "`"98499`"",
    "`"98374`"",
    "`"98391`"",
    "`"98375`"",
    "`"98270`"",
    "`"98271`"",
    "`"98223`"",
    "`"98516`"",
    "`"98503`"",
    "`"98584`"",
)

$startRow = 207
for ($i = 0; $i -lt $addresses.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "Walmart"
    $ws.Cells.Item($r, 2).Value = 2
    $ws.Cells.Item($r, 3).Value = $addresses[$i]
    $ws.Cells.Item($r, 4).Value = $zips[$i]
}

$ws.Range("D257").Select() | Out-Null